# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the last data row
# (row 5) on the zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-17 04:50:41"
$wsZhCn.Range("G5").Value = "2016-02-17 04:51:23"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-17 04:50:51"
$wsDeDe.Range("G5").Value = "2016-02-17 04:51:42"
